$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 12:22"

# Update country rows: reordering (Singapur, Marruecos, Republica del Chad moved up)
# and refreshed COVID-19 statistics for several countries.
$data = @{
    8 = @("Alemania", 153215, 86, 106800, 40840, 2908, 0, 5575)
    18 = @("Suiza", 28677, 181, 20600, 6526, 386, 2, 1551)
    24 = @("Austria", 15071, 69, 11872, 2677, 169, 0, 522)
    28 = @("Singapur", 12075, 897, 924, 11139, 26, 0, 12)
    29 = @("Chile", 11812, 0, 5804, 5840, 411, 0, 168)
    30 = @("Mexico", 11633, 1089, 2627, 7937, 378, 99, 1069)
    31 = @("Ecuador", 11183, 0, 1328, 9295, 141, 0, 560)
    35 = @("Rumania", 10417, 321, 2817, 7048, 241, 7, 552)
    55 = @("Marruecos", 3692, 124, 478, 3059, 1, 0, 155)
    56 = @("Luxemburgo", 3665, 0, 728, 2854, 27, 0, 83)
    62 = @("Kuwait", 2399, 0, 613, 1772, 55, 0, 14)
    69 = @("Uzbekistan", 1778, 20, 598, 1172, 8, 1, 8)
    96 = @("Libano", 696, 8, 140, 534, 46, 0, 22)
    161 = @("Libia", 60, 0, 18, 40, 0, 0, 2)
    169 = @("Republica del Chad", 40, 7, 8, 32, 0, 0, 0)
    170 = @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
    171 = @("Eritrea", 39, 0, 11, 28, 0, 0, 0)
    172 = @("San Martin (Parte Francesa)", 38, 0, 20, 16, 3, 0, 2)
    173 = @("Mongolia", 36, 1, 9, 27, 0, 0, 0)
    174 = @("Malaui", 33, 0, 3, 27, 1, 0, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
    $ws.Range("D$row").Value = $vals[3]
    $ws.Range("E$row").Value = $vals[4]
    $ws.Range("F$row").Value = $vals[5]
    $ws.Range("G$row").Value = $vals[6]
    $ws.Range("H$row").Value = $vals[7]
}
